$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Vscs")
$ws.Rows.Item(18).EntireRow.Insert()
Write-Host "A18:" $ws.Cells.Item(18,1).Value2
Write-Host "A19:" $ws.Cells.Item(19,1).Value2
Write-Host "A20:" $ws.Cells.Item(20,1).Value2
